$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("H2").Value = 'Americana O Klondike'
$ws.Range("I2").Value = 'Extra'
$ws.Range("J2").Value = 340
$ws.Range("O2").Value = 'Región de O''Higgins'

# Row 3
$ws.Range("D3").Value = Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("H3").Value = 'Americana O Klondike'
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 2000
$ws.Range("N3").Value = '$/unidad'
$ws.Range("O3").Value = 'Región de O''Higgins'
$ws.Range("P3").Value = 2000

# Row 4
$ws.Range("D4").Value = Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("H4").Value = 'Americana O Klondike'
$ws.Range("I4").Value = 'Segunda'
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 1500
$ws.Range("N4").Value = '$/unidad'
$ws.Range("O4").Value = 'Región de O''Higgins'
$ws.Range("P4").Value = 1500

# Row 5
$ws.Range("I5").Value = 'Tercera'
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 1000
$ws.Range("P5").Value = 1000

# Row 6
$ws.Range("D6").Value = Get-Date -Year 2021 -Month 10 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Range("H6").Value = 'Sin especificar'
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 800
$ws.Range("L6").Value = 800
$ws.Range("M6").Value = 800
$ws.Range("N6").Value = '$/kilo (volumen en unidades)'
$ws.Range("O6").Value = 'Perú'
$ws.Range("P6").Value = 800

# Row 7
$ws.Range("D7").Value = Get-Date -Year 2021 -Month 10 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 800
$ws.Range("L7").Value = 800
$ws.Range("M7").Value = 800
$ws.Range("N7").Value = '$/kilo (volumen en unidades)'
$ws.Range("O7").Value = 'Perú'
$ws.Range("P7").Value = 800

# Row 8
$ws.Range("D8").Value = Get-Date -Year 2021 -Month 10 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Range("H8").Value = 'Sin especificar'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = 800
$ws.Range("N8").Value = '$/kilo (volumen en unidades)'
$ws.Range("O8").Value = 'Perú'
$ws.Range("P8").Value = 800

# Row 9
$ws.Range("D9").Value = Get-Date -Year 2021 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("I9").Value = 'Extra'
$ws.Range("K9").Value = 2500
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2500
$ws.Range("P9").Value = 2500

# Row 10
$ws.Range("D10").Value = Get-Date -Year 2021 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 280
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 2000
$ws.Range("P10").Value = 2000

# Row 11
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = 5000
$ws.Range("P11").Value = 5000

# Row 12
$ws.Range("D12").Value = Get-Date -Year 2020 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("I12").Value = 'Segunda'
$ws.Range("J12").Value = 560
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("N12").Value = '$/unidad'
$ws.Range("O12").Value = 'Región de O''Higgins'
$ws.Range("P12").Value = 3000

# Row 13
$ws.Range("D13").Value = Get-Date -Year 2020 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("I13").Value = 'Tercera'
$ws.Range("J13").Value = 450
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 2000
$ws.Range("N13").Value = '$/unidad'
$ws.Range("O13").Value = 'Región de O''Higgins'
$ws.Range("P13").Value = 2000

# Row 14
$ws.Range("D14").Value = Get-Date -Year 2021 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = 800
$ws.Range("L14").Value = 800
$ws.Range("M14").Value = 800
$ws.Range("N14").Value = '$/kilo (volumen en unidades)'
$ws.Range("O14").Value = 'Perú'
$ws.Range("P14").Value = 800

# Row 15
$ws.Range("D15").Value = Get-Date -Year 2021 -Month 11 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("K15").Value = 800
$ws.Range("L15").Value = 800
$ws.Range("M15").Value = 800
$ws.Range("N15").Value = '$/kilo (volumen en unidades)'
$ws.Range("O15").Value = 'Perú'
$ws.Range("P15").Value = 800

# Row 16
$ws.Range("D16").Value = Get-Date -Year 2021 -Month 4 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("J16").Value = 180

# Row 17
$ws.Range("D17").Value = Get-Date -Year 2021 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("J17").Value = 150

# Row 18
$ws.Range("D18").Value = Get-Date -Year 2021 -Month 4 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 2500
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = 2500
$ws.Range("N18").Value = '$/unidad'
$ws.Range("P18").Value = 2500

# Row 19
$ws.Range("D19").Value = Get-Date -Year 2021 -Month 10 -Day 26 -Hour 0 -Minute 0 -Second 0

# Row 20
$ws.Range("D20").Value = Get-Date -Year 2021 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = 800
$ws.Range("L20").Value = 800
$ws.Range("M20").Value = 800
$ws.Range("N20").Value = '$/kilo (volumen en unidades)'
$ws.Range("O20").Value = 'Perú'
$ws.Range("P20").Value = 800

# Row 21
$ws.Range("D21").Value = Get-Date -Year 2020 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("I21").Value = 'Extra'
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 3500
$ws.Range("L21").Value = 3500
$ws.Range("M21").Value = 3500
$ws.Range("P21").Value = 3500

# Row 22
$ws.Range("D22").Value = Get-Date -Year 2020 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 3000
$ws.Range("N22").Value = '$/unidad'
$ws.Range("O22").Value = 'Región de O''Higgins'
$ws.Range("P22").Value = 3000
